# WorkingHours.xlsx update:
#   - Each person's sheet gets a new "working hours" entry for 2024-04-01
#     (serial 45383): PaulSchein & ReneMifka worked "3h" on
#     "Themes Router und http test Requests"; LukasPerger worked "1h" on
#     the same task.
#   - The new date cell is formatted like the existing date column (copy
#     the format from the row above so it reuses the same date style).
#   - LukasPerger's sheet (the 3rd tab) ends up the active/selected sheet,
#     matching the new selections left on each sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: PaulSchein ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A5").Copy() | Out-Null
$ws1.Range("A6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws1.Range("A6").Value = 45383
$ws1.Range("B6").Value = "3h"
$ws1.Range("C6").Value = "Themes Router und http test Requests"
$ws1.Range("C6").Select() | Out-Null

# --- Sheet 2: ReneMifka ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A5").Copy() | Out-Null
$ws2.Range("A6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws2.Range("A6").Value = 45383
$ws2.Range("B6").Value = "3h"
$ws2.Range("C6").Value = "Themes Router und http test Requests"
$ws2.Range("C7").Select() | Out-Null

# --- Sheet 3: LukasPerger ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A6").Copy() | Out-Null
$ws3.Range("A7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws3.Range("A7").Value = 45383
$ws3.Range("B7").Value = "1h"
$ws3.Range("C7").Value = "Themes Router und http test Requests"
$ws3.Range("A8").Select() | Out-Null
$ws3.Activate() | Out-Null
